$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = "Normal"
}

Set-TextValue "D2" "314.68"
Set-TextValue "E2" "3.10%"
Set-TextValue "D3" "35.50"
Set-TextValue "E3" "-1.95%"
Set-TextValue "D4" "5.086"
Set-TextValue "E4" "1.31%"
Set-TextValue "D5" "0.08183"
Set-TextValue "E5" "3.87%"
Set-TextValue "D6" "2.106"
Set-TextValue "E6" "-1.09%"
Set-TextValue "D7" "7.971"
Set-TextValue "E7" "0.30%"
Set-TextValue "D8" "4.139"
Set-TextValue "E8" "0.21%"
Set-TextValue "E9" "0.64%"
Set-TextValue "E10" "6.69%"
Set-TextValue "D11" "0.1935"
Set-TextValue "E11" "4.78%"
Set-TextValue "D12" "0.09199"
Set-TextValue "E12" "6.27%"
Set-TextValue "D13" "0.03635"
Set-TextValue "E13" "2.15%"
Set-TextValue "D14" "0.09893"
Set-TextValue "E14" "-0.18%"
Set-TextValue "D15" "0.001435"
Set-TextValue "E15" "-0.33%"
Set-TextValue "D16" "0.005656"
Set-TextValue "E16" "-0.73%"
Set-TextValue "D17" "3.475"
Set-TextValue "E17" "0.12%"
Set-TextValue "D18" "2.976"
Set-TextValue "E18" "8.12%"
Set-TextValue "D19" "0.3413"
Set-TextValue "E19" "1.14%"
Set-TextValue "D20" "0.1301"
Set-TextValue "E20" "-3.56%"
Set-TextValue "D21" "5.104"
Set-TextValue "E21" "-0.99%"
Set-TextValue "E22" "0.22%"
Set-TextValue "D23" "0.04547"
Set-TextValue "E23" "-0.79%"
Set-TextValue "D24" "0.001232"
Set-TextValue "E24" "0.05%"
Set-TextValue "D25" "0.004793"
Set-TextValue "E25" "-0.27%"
Set-TextValue "D26" "0.0001251"
Set-TextValue "E26" "-3.78%"
Set-TextValue "D27" "0.0004452"
Set-TextValue "D39" "0.01999"
Set-TextValue "E39" "7.85%"
Set-TextValue "D40" "0.04932"
Set-TextValue "D41" "0.007544"
Set-TextValue "E41" "-3.23%"
Set-TextValue "D42" "0.1385"
Set-TextValue "E42" "0.16%"
Set-TextValue "D43" "0.007941"
Set-TextValue "E43" "2.67%"
Set-TextValue "D44" "0.002222"
Set-TextValue "E44" "2.76%"
Set-TextValue "D45" "0.01157"
Set-TextValue "E45" "1.86%"
Set-TextValue "D46" "0.00006596"
Set-TextValue "E46" "3.46%"
Set-TextValue "E47" "0.06%"
Set-TextValue "D48" "185.53"
Set-TextValue "E48" "257.28%"
Set-TextValue "E49" "-10.50%"
Set-TextValue "D50" "0.00002102"
Set-TextValue "E50" "0.06%"
Set-TextValue "D51" "0.0002002"
Set-TextValue "E51" "0.06%"
